# Update "想去人数" (F column) figures that were refreshed by the scraper.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 0
$ws1.Range("F3").Value  = 0
$ws1.Range("F4").Value  = 0
$ws1.Range("F6").Value  = 0
$ws1.Range("F7").Value  = 6754
$ws1.Range("F9").Value  = 0
$ws1.Range("F10").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("F18").Value = 26
$ws1.Range("F20").Value = 0
$ws1.Range("F21").Value = 0
$ws1.Range("F22").Value = 0
$ws1.Range("F23").Value = 549
$ws1.Range("F25").Value = 218

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 0

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 0
$ws4.Range("F4").Value  = 0
$ws4.Range("F7").Value  = 6754
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 0
$ws4.Range("F24").Value = 0
$ws4.Range("F25").Value = 549
$ws4.Range("F27").Value = 0
